$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.413.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.43%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.73%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0632"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "

# Row 12
$ws.Range("E12").Value = "  -0.71%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.670.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.351.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

# Row 23
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0511"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "

# Row 32
$ws.Range("E32").Value = "  +0.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "

# Row 34
$ws.Range("E34").Value = "  +1.45%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.258.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.47%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0179"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.846"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "

# Row 41
$ws.Range("E41").Value = "  -1.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "

# Row 43
$ws.Range("E43").Value = "  +4.75%  "

# Row 44
$ws.Range("E44").Value = "  -0.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0514"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "

# Row 51
$ws.Range("E51").Value = "  -0.77%  "

Write-Host "Applied cryptos list update"